# Update the "Return_with_prediction" (G), "return_pct_change" (H) and the
# single changed "mean_return_pct_change" (I2) values in Sheet1, rows 2-29,
# to reflect the recomputed figures from the latest recurrence run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.04281337908013565
$ws.Range("G3").Value = 0.05765081963801109
$ws.Range("G4").Value = -0.4535044621812089
$ws.Range("G5").Value = -0.4868696565836673
$ws.Range("G6").Value = 0.232143033415955
$ws.Range("G7").Value = 0.2299493561456461
$ws.Range("G8").Value = 0.1641386126071218
$ws.Range("G9").Value = 0.1736065283855741
$ws.Range("G10").Value = -0.01233644169476444
$ws.Range("G11").Value = -0.03227087510197125
$ws.Range("G12").Value = 0.1316122180289581
$ws.Range("G13").Value = 0.1330185301764849
$ws.Range("G14").Value = 0.2560563878134899
$ws.Range("G15").Value = 0.2574673972360753
$ws.Range("G16").Value = 0.1328170162102214
$ws.Range("G17").Value = 0.1410961765994528
$ws.Range("G18").Value = -0.02627778221040918
$ws.Range("G19").Value = -0.01420724301156604
$ws.Range("G20").Value = 0.1370448739116555
$ws.Range("G21").Value = 0.1426524975884551
$ws.Range("G22").Value = 0.1730249879388232
$ws.Range("G23").Value = 0.1773148304219258
$ws.Range("G24").Value = -0.1091503220107316
$ws.Range("G25").Value = -0.1032928758646516
$ws.Range("G26").Value = 0.2322117338454718
$ws.Range("G27").Value = 0.2304701074549608
$ws.Range("G28").Value = 0.07395570645338267
$ws.Range("G29").Value = 0.06785823106279887
$ws.Range("H2").Value = -11.23195571966703
$ws.Range("H3").Value = 50.30220202550939
$ws.Range("H4").Value = -0.1626602872602654
$ws.Range("H5").Value = -1.677963871920282
$ws.Range("H6").Value = -0.6371642686129091
$ws.Range("H7").Value = 4.250146675925488
$ws.Range("H8").Value = -1.599354075335968
$ws.Range("H9").Value = 0.9273799772762804
$ws.Range("H10").Value = -160.2909367988167
$ws.Range("H11").Value = -120.1994133928426
$ws.Range("H12").Value = -3.73754529118185
$ws.Range("H13").Value = 6.72425434141914
$ws.Range("H14").Value = 3.526781882335527
$ws.Range("H15").Value = 1.895854033329996
$ws.Range("H16").Value = -13.45693154753818
$ws.Range("H17").Value = -6.568227508947857
$ws.Range("H18").Value = -60.47592457189541
$ws.Range("H19").Value = -1587.192954197784
$ws.Range("H20").Value = -1.16635917747376
$ws.Range("H21").Value = -0.3088711863501672
$ws.Range("H22").Value = -7.084407275830682
$ws.Range("H23").Value = -1.195624527095493
$ws.Range("H24").Value = -15.62869399526264
$ws.Range("H25").Value = -3.698904186270948
$ws.Range("H26").Value = 0.9000118365330361
$ws.Range("H27").Value = -0.9038614400259046
$ws.Range("H28").Value = 25.76803239806168
$ws.Range("H29").Value = -3.86244787534923
$ws.Range("I2").Value = -68.0994835008954
